$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing "Regression" header in column G (row 1)
$ws.Range("G1").Value = "Regression"

# Update the active selection to reflect where the user ended up (G2)
$ws.Range("G2").Select()
